# Apply updated crypto price/volume figures (symbol-list refresh).
# Each D/E cell in the source sheet is stored as text (inlineStr), so we
# write the new value with a leading apostrophe to force Excel to keep it
# as text instead of auto-converting to a number/percentage, then reset
# the cell style back to "Normal" so no stray number-format style sticks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "315.05"
Set-TextValue "E2" "2.45%"
Set-TextValue "D3" "40.77"
Set-TextValue "E3" "-0.73%"
Set-TextValue "D4" "5.169"
Set-TextValue "E4" "-1.45%"
Set-TextValue "D5" "0.07590"
Set-TextValue "E5" "-0.97%"
Set-TextValue "D6" "4.327"
Set-TextValue "E6" "0.22%"
Set-TextValue "D7" "1.666"
Set-TextValue "E7" "2.46%"
Set-TextValue "D8" "0.9274"
Set-TextValue "E8" "1.00%"
Set-TextValue "E10" "-2.39%"
Set-TextValue "D11" "0.1819"
Set-TextValue "E11" "-0.35%"
Set-TextValue "D12" "0.09044"
Set-TextValue "E12" "-0.42%"
Set-TextValue "D13" "0.04153"
Set-TextValue "E13" "-2.43%"
Set-TextValue "E14" "0.24%"
Set-TextValue "D15" "0.001291"
Set-TextValue "E15" "2.41%"
Set-TextValue "E16" "0.28%"
Set-TextValue "D18" "3.331"
Set-TextValue "E18" "-0.68%"
Set-TextValue "E19" "0.66%"
Set-TextValue "D20" "7.568"
Set-TextValue "E20" "3.46%"
Set-TextValue "E21" "-2.37%"
Set-TextValue "E22" "-3.04%"
Set-TextValue "D23" "0.04032"
Set-TextValue "D24" "0.001273"
Set-TextValue "E24" "0.78%"
Set-TextValue "D25" "0.004051"
Set-TextValue "E25" "-7.08%"
Set-TextValue "D26" "0.0001270"
Set-TextValue "E26" "-0.22%"
Set-TextValue "D38" "0.02411"
Set-TextValue "E38" "-2.46%"
Set-TextValue "D39" "0.05164"
Set-TextValue "E39" "-2.27%"
Set-TextValue "D40" "0.007722"
Set-TextValue "E40" "-1.57%"
Set-TextValue "D41" "0.1301"
Set-TextValue "E41" "-0.92%"
Set-TextValue "D42" "0.007607"
Set-TextValue "E42" "16.93%"
Set-TextValue "D43" "0.003300"
Set-TextValue "E43" "72.41%"
Set-TextValue "D44" "0.008555"
Set-TextValue "E44" "10.73%"
Set-TextValue "D45" "0.3407"
Set-TextValue "E45" "11.32%"
Set-TextValue "D46" "0.00006583"
Set-TextValue "E46" "-1.97%"
Set-TextValue "E47" "-0.21%"
Set-TextValue "D48" "0.2752"
Set-TextValue "E48" "62.03%"
Set-TextValue "D49" "0.004203"
Set-TextValue "E49" "2.51%"
Set-TextValue "E50" "-0.21%"
Set-TextValue "E51" "-0.21%"
